$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 57253.75
$ws.Range("J26").Value = 57253.75
$ws.Range("L26").Value = 57253.75
$ws.Range("N26").Value = -57941.75

$ws.Range("H62").Value = 9006.806
$ws.Range("I62").Value = 11112.926
$ws.Range("J62").Value = 2688.4443
$ws.Range("K62").Value = 11112.926
$ws.Range("L62").Value = 2688.4443
$ws.Range("M62").Value = -10488.926
$ws.Range("N62").Value = -3936.4443

$ws.Range("H65").Value = 9006.806
$ws.Range("I65").Value = 11112.926
$ws.Range("J65").Value = 2688.4443
$ws.Range("K65").Value = 55564.63
$ws.Range("L65").Value = 13442.2215
$ws.Range("M65").Value = -52444.63
$ws.Range("N65").Value = -19682.2215

$ws.Range("H98").Value = 1513.762
$ws.Range("I98").Value = 930.1875
$ws.Range("J98").Value = 3381.2
$ws.Range("K98").Value = 930.1875
$ws.Range("L98").Value = 3381.2
$ws.Range("M98").Value = 567.8125
$ws.Range("N98").Value = -6377.2

$ws.Range("H122").Value = 1513.762
$ws.Range("I122").Value = 930.1875
$ws.Range("J122").Value = 3381.2
$ws.Range("K122").Value = 2790.5625
$ws.Range("L122").Value = 10143.6
$ws.Range("M122").Value = -340.5625
$ws.Range("N122").Value = -15043.6

$ws.Range("H138").Value = 1699.14
$ws.Range("I138").Value = 963.0714
$ws.Range("J138").Value = 1985.3889
$ws.Range("K138").Value = 2889.2142
$ws.Range("L138").Value = 5956.1667
$ws.Range("M138").Value = 2250.7858
$ws.Range("N138").Value = -16236.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1082.0667
$ws.Range("I2").Value = 1102.5834
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1102.5834
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -989.5834
$ws.Range("N2").Value = -1226

$ws.Range("H52").Value = 39780
$ws.Range("J52").Value = 39780
$ws.Range("L52").Value = 39780
$ws.Range("N52").Value = -40416

$ws.Range("H116").Value = 1082.0667
$ws.Range("I116").Value = 1102.5834
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 1102.5834
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 1191.4166
$ws.Range("N116").Value = -5588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1082.0667
$ws.Range("I3").Value = 1102.5834
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 1102.5834
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -988.5834
$ws.Range("N3").Value = -1228

$ws.Range("H94").Value = 1579.1875
$ws.Range("I94").Value = 480.7
$ws.Range("J94").Value = 3410
$ws.Range("K94").Value = 480.7
$ws.Range("L94").Value = 3410
$ws.Range("M94").Value = -29.69999999999999
$ws.Range("N94").Value = -4312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24444.66
$ws.Range("I31").Value = 29163.486
$ws.Range("J31").Value = 6985
$ws.Range("K31").Value = 29163.486
$ws.Range("L31").Value = 6985
$ws.Range("M31").Value = -28868.486
$ws.Range("N31").Value = -7575

$ws.Range("H34").Value = 24444.66
$ws.Range("I34").Value = 29163.486
$ws.Range("J34").Value = 6985
$ws.Range("K34").Value = 29163.486
$ws.Range("L34").Value = 6985
$ws.Range("M34").Value = -28961.486
$ws.Range("N34").Value = -7389

$ws.Range("H35").Value = 16953
$ws.Range("I35").Value = 922.7778
$ws.Range("J35").Value = 53021
$ws.Range("K35").Value = 922.7778
$ws.Range("L35").Value = 53021
$ws.Range("M35").Value = -628.7778
$ws.Range("N35").Value = -53609

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
# clear M105
$ws.Range("M105").ClearContents()
# clear N105
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 13444.083
$ws.Range("I5").Value = 6096.4
$ws.Range("J5").Value = 50182.5
$ws.Range("K5").Value = 18289.2
$ws.Range("L5").Value = 150547.5
$ws.Range("M5").Value = -18177.2
$ws.Range("N5").Value = -150771.5

$ws.Range("H24").Value = 2627.4285
$ws.Range("J24").Value = 2627.4285
$ws.Range("L24").Value = 7882.2855
$ws.Range("N24").Value = -8342.2855

$ws.Range("H122").Value = 368.13333
$ws.Range("I122").Value = 251.64285
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 2264.78565
$ws.Range("L122").Value = 17991
$ws.Range("M122").Value = 185.2143499999997
$ws.Range("N122").Value = -22891

$ws.Range("H135").Value = 13444.083
$ws.Range("I135").Value = 6096.4
$ws.Range("J135").Value = 50182.5
$ws.Range("K135").Value = 54867.6
$ws.Range("L135").Value = 451642.5
$ws.Range("M135").Value = -52332.6
$ws.Range("N135").Value = -456712.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 61666.668
$ws.Range("J88").Value = 61666.668
$ws.Range("L88").Value = 61666.668
$ws.Range("N88").Value = -62568.668

$ws.Range("H91").Value = 61666.668
$ws.Range("J91").Value = 61666.668
$ws.Range("L91").Value = 61666.668
$ws.Range("N91").Value = -64786.668

$ws.Range("H141").Value = 41466
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 44532.57
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 44532.57
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -54892.57

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 912498.5600000001
$ws.Range("I136").Value = 2003012.8
$ws.Range("J136").Value = 3736.6667
$ws.Range("K136").Value = 6009038.4
$ws.Range("L136").Value = 11210.0001
$ws.Range("M136").Value = -6006488.4
$ws.Range("N136").Value = -16310.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 70017
$ws.Range("J21").Value = 70017
$ws.Range("L21").Value = 70017
$ws.Range("N21").Value = -70487

$ws.Range("H35").Value = 70017
$ws.Range("J35").Value = 70017
$ws.Range("L35").Value = 70017
$ws.Range("N35").Value = -70597

$ws.Range("H140").Value = 42500
$ws.Range("I140").Value = 20000
$ws.Range("J140").Value = 50000
$ws.Range("K140").Value = 20000
$ws.Range("L140").Value = 50000
$ws.Range("M140").Value = -14820
$ws.Range("N140").Value = -60360

$ws.Range("H141").Value = 54750
$ws.Range("I141").Value = 40000
$ws.Range("J141").Value = 56857.145
$ws.Range("K141").Value = 40000
$ws.Range("L141").Value = 56857.145
$ws.Range("M141").Value = -34820
$ws.Range("N141").Value = -67217.14499999999
